# The document contains one table; its data rows (1, 5, 9, 13, 17 in the
# Word 1-based Rows index) each hold 5 division-expression cells. We update
# each cell's text in place by row/column coordinates so that duplicate
# expressions (e.g. the two "26÷7=" cells) are replaced independently with
# their own distinct target values, matching the diff exactly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "99÷6="
$t.Cell(1,2).Range.Text  = "12÷2="
$t.Cell(1,3).Range.Text  = "16÷6="
$t.Cell(1,4).Range.Text  = "99÷4="
$t.Cell(1,5).Range.Text  = "38÷2="

$t.Cell(5,1).Range.Text  = "72÷4="
$t.Cell(5,2).Range.Text  = "10÷7="
$t.Cell(5,3).Range.Text  = "51÷3="
$t.Cell(5,4).Range.Text  = "84÷5="
$t.Cell(5,5).Range.Text  = "56÷8="

$t.Cell(9,1).Range.Text  = "17÷9="
$t.Cell(9,2).Range.Text  = "87÷8="
$t.Cell(9,3).Range.Text  = "51÷8="
$t.Cell(9,4).Range.Text  = "25÷4="
$t.Cell(9,5).Range.Text  = "87÷5="

$t.Cell(13,1).Range.Text = "54÷2="
$t.Cell(13,2).Range.Text = "23÷8="
$t.Cell(13,3).Range.Text = "33÷8="
$t.Cell(13,4).Range.Text = "46÷7="
$t.Cell(13,5).Range.Text = "44÷4="

$t.Cell(17,1).Range.Text = "29÷6="
$t.Cell(17,2).Range.Text = "83÷5="
$t.Cell(17,3).Range.Text = "81÷3="
$t.Cell(17,4).Range.Text = "57÷7="
$t.Cell(17,5).Range.Text = "30÷7="
